$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.311.20"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.490.60"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'321.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'108.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'39.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'18.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "'7.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "2.878.08"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "2.495.35"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "'0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "47.236.64"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "'13.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("D20").Value = "'6.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'2.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.39%  "
$ws.Range("D23").Value = "'70.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'247.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "'9.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'34.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "'49.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "'20.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "'0.0783"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "'23.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.80%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").Value = "'119.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").Value = "'0.0297"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "1.999.35"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").Value = "'3.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'2.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "'5.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").Value = "'56.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.07%  "
